$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row4_C = New-Object 'object[,]' 1,3
$row4_C[0,0] = 0.2886288628862886
$row4_C[0,1] = 0.9879147914791477
$row4_C[0,2] = -0.9999879987998799
$ws.Range("C4:E4").Value = $row4_C

$row4_G = New-Object 'object[,]' 1,7
$row4_G[0,0] = -0.7213921392139213
$row4_G[0,1] = 0.3545274527452745
$row4_G[0,2] = -0.01120912091209121
$row4_G[0,3] = -0.003924392439243924
$row4_G[0,4] = -0.06499849984998499
$row4_G[0,5] = -0.05506150615061506
$row4_G[0,6] = -0.2575937593759375
$ws.Range("G4:M4").Value = $row4_G

$row4_O = New-Object 'object[,]' 1,2
$row4_O[0,0] = -0.07167116711671166
$row4_O[0,1] = 0.03185118511851185
$ws.Range("O4:P4").Value = $row4_O

$row5_C = New-Object 'object[,]' 1,3
$row5_C[0,0] = -0.127980798079808
$row5_C[0,1] = -0.1746654665466546
$row5_C[0,2] = 0.1735733573357336
$ws.Range("C5:E5").Value = $row5_C

$row5_G = New-Object 'object[,]' 1,7
$row5_G[0,0] = 0.1646324632463246
$row5_G[0,1] = -0.148982898289829
$row5_G[0,2] = -0.1615121512151215
$row5_G[0,3] = 0.04241224122412241
$row5_G[0,4] = 0.1749174917491749
$row5_G[0,5] = -0.01184518451845184
$row5_G[0,6] = 0.1174677467746775
$ws.Range("G5:M5").Value = $row5_G

$row5_O = New-Object 'object[,]' 1,2
$row5_O[0,0] = 0.1235043504350435
$row5_O[0,1] = -0.1793939393939394
$ws.Range("O5:P5").Value = $row5_O

$row6_C = New-Object 'object[,]' 1,3
$row6_C[0,0] = 0.1322532253225322
$row6_C[0,1] = 0.1131953195319532
$row6_C[0,2] = -0.1153315331533153
$ws.Range("C6:E6").Value = $row6_C

$row6_G = New-Object 'object[,]' 1,7
$row6_G[0,0] = -0.08160816081608159
$row6_G[0,1] = 0.04373237323732372
$row6_G[0,2] = -0.05526552655265525
$row6_G[0,3] = -0.001188118811881188
$row6_G[0,4] = 0.0463006300630063
$row6_G[0,5] = 0.0008880888088808879
$row6_G[0,6] = 0.07713171317131713
$ws.Range("G6:M6").Value = $row6_G

$row6_O = New-Object 'object[,]' 1,2
$row6_O[0,0] = 0.1208160816081608
$row6_O[0,1] = -0.01336933693369337
$ws.Range("O6:P6").Value = $row6_O

$row7_C = New-Object 'object[,]' 1,3
$row7_C[0,0] = 0.1062466246624662
$row7_C[0,1] = 0.09689768976897688
$row7_C[0,2] = 0.04159615961596159
$ws.Range("C7:E7").Value = $row7_C

$row7_G = New-Object 'object[,]' 1,7
$row7_G[0,0] = 0.4447884788478847
$row7_G[0,1] = 0.008832883288328831
$row7_G[0,2] = 0.004116411641164116
$row7_G[0,3] = -0.1628082808280828
$row7_G[0,4] = 0.2631623162316231
$row7_G[0,5] = 0.9896549654965496
$row7_G[0,6] = -0.1831743174317431
$ws.Range("G7:M7").Value = $row7_G

$row7_O = New-Object 'object[,]' 1,2
$row7_O[0,0] = 0.08698469846984698
$row7_O[0,1] = 0.2030723072307231
$ws.Range("O7:P7").Value = $row7_O

$row8_C = New-Object 'object[,]' 1,3
$row8_C[0,0] = 0.07792379237923792
$row8_C[0,1] = -0.05966996699669967
$row8_C[0,2] = 0.06756675667566756
$ws.Range("C8:E8").Value = $row8_C

$row8_G = New-Object 'object[,]' 1,7
$row8_G[0,0] = -0.4363756375637562
$row8_G[0,1] = 0.9027662766276627
$row8_G[0,2] = 0.7969636963696368
$row8_G[0,3] = -0.02419441944194419
$row8_G[0,4] = 0.05653765376537653
$row8_G[0,5] = 0.008412841284128413
$row8_G[0,6] = 0.1164956495649565
$ws.Range("G8:M8").Value = $row8_G

$row8_O = New-Object 'object[,]' 1,2
$row8_O[0,0] = 0.0293069306930693
$row8_O[0,1] = 0.0686108610861086
$ws.Range("O8:P8").Value = $row8_O

$row9_C = New-Object 'object[,]' 1,3
$row9_C[0,0] = 0.7948394839483948
$row9_C[0,1] = 0.03437143714371436
$row9_C[0,2] = -0.04020402040204019
$ws.Range("C9:E9").Value = $row9_C

$row9_G = New-Object 'object[,]' 1,7
$row9_G[0,0] = -0.04452445244524452
$row9_G[0,1] = 0.007908790879087909
$row9_G[0,2] = 0.1024542454245424
$row9_G[0,3] = -0.2369636963696369
$row9_G[0,4] = -0.4938133813381337
$row9_G[0,5] = -0.03239123912391239
$row9_G[0,6] = -0.0266066606660666
$ws.Range("G9:M9").Value = $row9_G

$row9_O = New-Object 'object[,]' 1,2
$row9_O[0,0] = -0.04922892289228922
$row9_O[0,1] = -0.04145214521452145
$ws.Range("O9:P9").Value = $row9_O

$row10_C = New-Object 'object[,]' 1,3
$row10_C[0,0] = 0.2643144314431443
$row10_C[0,1] = -0.05465346534653465
$row10_C[0,2] = 0.06719471947194719
$ws.Range("C10:E10").Value = $row10_C

$row10_G = New-Object 'object[,]' 1,7
$row10_G[0,0] = 0.03553555355535553
$row10_G[0,1] = 0.01431743174317432
$row10_G[0,2] = -0.05666966696669666
$row10_G[0,3] = 0.01874587458745874
$row10_G[0,4] = 0.8615181518151815
$row10_G[0,5] = 0.07876387638763875
$row10_G[0,6] = 0.0341074107410741
$ws.Range("G10:M10").Value = $row10_G

$row10_O = New-Object 'object[,]' 1,2
$row10_O[0,0] = 0.03491149114911491
$row10_O[0,1] = -0.01278127812781278
$ws.Range("O10:P10").Value = $row10_O

$row11_C = New-Object 'object[,]' 1,3
$row11_C[0,0] = -0.1255925592559256
$row11_C[0,1] = 0.08823282328232822
$row11_C[0,2] = -0.1003300330033003
$ws.Range("C11:E11").Value = $row11_C

$row11_G = New-Object 'object[,]' 1,7
$row11_G[0,0] = -0.01083708370837084
$row11_G[0,1] = -0.1441224122412241
$row11_G[0,2] = -0.1155475547554755
$row11_G[0,3] = 0.1300690069006901
$row11_G[0,4] = -0.07433543354335433
$row11_G[0,5] = -0.1348454845484548
$row11_G[0,6] = 0.0381998199819982
$ws.Range("G11:M11").Value = $row11_G

$row11_O = New-Object 'object[,]' 1,2
$row11_O[0,0] = -0.02364236423642364
$row11_O[0,1] = -0.04432043204320432
$ws.Range("O11:P11").Value = $row11_O

$row12_C = New-Object 'object[,]' 1,3
$row12_C[0,0] = -0.1022982298229823
$row12_C[0,1] = -0.1913591359135913
$row12_C[0,2] = 0.1784338433843384
$ws.Range("C12:E12").Value = $row12_C

$row12_G = New-Object 'object[,]' 1,7
$row12_G[0,0] = 0.1765736573657365
$row12_G[0,1] = -0.1659525952595259
$row12_G[0,2] = -0.1073267326732673
$row12_G[0,3] = -0.04138013801380137
$row12_G[0,4] = 0.01261326132613261
$row12_G[0,5] = 0.02282628262826283
$row12_G[0,6] = 0.01850585058505851
$ws.Range("G12:M12").Value = $row12_G

$row12_O = New-Object 'object[,]' 1,2
$row12_O[0,0] = 0.01154515451545154
$row12_O[0,1] = 0.03655565556555655
$ws.Range("O12:P12").Value = $row12_O

$row13_C = New-Object 'object[,]' 1,3
$row13_C[0,0] = 0.2552295229522952
$row13_C[0,1] = 0.07637563756375637
$row13_C[0,2] = -0.08162016201620162
$ws.Range("C13:E13").Value = $row13_C

$row13_G = New-Object 'object[,]' 1,7
$row13_G[0,0] = -0.03401140114011401
$row13_G[0,1] = 0.01492949294929493
$row13_G[0,2] = 0.564128412841284
$row13_G[0,3] = -0.00252025202520252
$row13_G[0,4] = -0.2173897389738974
$row13_G[0,5] = -0.01333333333333333
$row13_G[0,6] = 0.0572097209720972
$ws.Range("G13:M13").Value = $row13_G

$row13_O = New-Object 'object[,]' 1,2
$row13_O[0,0] = -0.06673867386738673
$row13_O[0,1] = -0.04013201320132012
$ws.Range("O13:P13").Value = $row13_O

$row14_C = New-Object 'object[,]' 1,3
$row14_C[0,0] = -0.2215541554155415
$row14_C[0,1] = -0.01572157215721572
$row14_C[0,2] = -0.00288028802880288
$ws.Range("C14:E14").Value = $row14_C

$row14_G = New-Object 'object[,]' 1,7
$row14_G[0,0] = -0.08727272727272728
$row14_G[0,1] = 0.02412241224122412
$row14_G[0,2] = -0.05876987698769877
$row14_G[0,3] = 0.1085628562856285
$row14_G[0,4] = -0.139021902190219
$row14_G[0,5] = -0.1684728472847285
$row14_G[0,6] = 0.09477347734773477
$ws.Range("G14:M14").Value = $row14_G

$row14_O = New-Object 'object[,]' 1,2
$row14_O[0,0] = -0.02085808580858085
$row14_O[0,1] = -0.02761476147614761
$ws.Range("O14:P14").Value = $row14_O

Write-Output "done"